$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# This worksheet is a table of simulated "HKL" entries. The author reran
# the simulation notebook which (a) renamed "Thomas Hex" to "Matthies Hex"
# and (b) added two new entries, "Holden" and "Rizzie Spiral", right after
# "Spiral5" (i.e. as rows 4 and 5), pushing all following rows down by two.
# -----------------------------------------------------------------------

# 1) Insert two new blank rows at row 4, pushing the existing rows 4-29
#    down to rows 6-31 (carrying all of their data/formatting with them).
$ws.Range("A4:A5").EntireRow.Insert()

# 2) Copy formatting from the (now shifted) rows below into the two new
#    blank rows so they match the rest of the table's look (bordered,
#    bold, centered label column, plain numeric columns).
$ws.Range("A6:W6").Copy()
$ws.Range("A4:W4").PasteSpecial(-4122)
$ws.Range("A7:W7").Copy()
$ws.Range("A5:W5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Fill in the label columns for the two new rows.
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"

# 4) Fill in the freshly-simulated numeric data for the two new rows.
$row4 = New-Object 'object[,]' 1,21
$row4[0,0]  = 1.008988713663024
$row4[0,1]  = 0.9704984160446846
$row4[0,2]  = 1.016818949156604
$row4[0,3]  = 0.9704984160446846
$row4[0,4]  = 0.9951246855042575
$row4[0,5]  = 1.006481524973545
$row4[0,6]  = 0.981090822293042
$row4[0,7]  = 1.016818949156604
$row4[0,8]  = 1.016818949156604
$row4[0,9]  = 0.9937510412163307
$row4[0,10] = 1.011125095538023
$row4[0,11] = 1.016818949156604
$row4[0,12] = 1.008988713663024
$row4[0,13] = 0.9897435648538542
$row4[0,14] = 1.001369877439677
$row4[0,15] = 0.998768692954771
$row4[0,16] = 0.9910793903080131
$row4[0,17] = 0.998768692954771
$row4[0,18] = 0.997514280020161
$row4[0,19] = 1.00137521384745
$row4[0,20] = 0.9979849060486887
$ws.Range("C4:W4").Value = $row4

$row5 = New-Object 'object[,]' 1,21
$row5[0,0]  = 1.033474386883205
$row5[0,1]  = 0.870239079577972
$row5[0,2]  = 1.085745272975651
$row5[0,3]  = 0.870239079577972
$row5[0,4]  = 0.9742911369010427
$row5[0,5]  = 1.032000259982808
$row5[0,6]  = 0.9154509069023515
$row5[0,7]  = 1.085745272975651
$row5[0,8]  = 1.085745272975651
$row5[0,9]  = 0.9718044450903923
$row5[0,10] = 1.048553664528656
$row5[0,11] = 1.085745272975651
$row5[0,12] = 1.033474386883205
$row5[0,13] = 0.9518567332305883
$row5[0,14] = 1.002639415986798
$row5[0,15] = 0.9964862464789425
$row5[0,16] = 0.9585059705171896
$row5[0,17] = 0.9964862464789425
$row5[0,18] = 0.9903157961318049
$row5[0,19] = 1.009401691500574
$row5[0,20] = 0.9914448941052598
$ws.Range("C5:W5").Value = $row5

# 5) Rename "Thomas Hex" -> "Matthies Hex" wherever it occurs in the sheet.
$found = $ws.Cells.Find("Thomas Hex")
if ($found) {
    $found.Value = "Matthies Hex"
}

# 6) Column A is a simple running index (0, 1, 2, ...) independent of the
#    row-shift that happened above; re-number it sequentially for every
#    data row now that the table has grown from 28 to 30 entries.
$idx = New-Object 'object[,]' 30,1
for ($i = 0; $i -lt 30; $i++) {
    $idx[$i,0] = $i
}
$ws.Range("A2:A31").Value = $idx

# 7) Make sure the sheet's declared dimension covers the new rows.
$ws.Range("A1:W31").Value = $ws.Range("A1:W31").Value()
